$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data held in columns A, B, D, E, F, G, H, J, Q, R between rows 5 and 6.
$cols = @("A", "B", "D", "E", "F", "G", "H", "J", "Q", "R")

foreach ($col in $cols) {
    $addr5 = "{0}5" -f $col
    $addr6 = "{0}6" -f $col
    $v5 = $ws.Range($addr5).Value2
    $v6 = $ws.Range($addr6).Value2
    $ws.Range($addr5).Value2 = $v6
    $ws.Range($addr6).Value2 = $v5
}
